$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting so numeric-looking price strings
# (e.g. "1.00", "215.60") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.840.28"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "1.666.89"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "215.60"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("E6").Value = "  +5.34%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "20.21"
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("D11").Value = "0.0890"
$ws.Range("E11").Value = "  +3.91%  "
$ws.Range("D12").Value = "1.904.43"
$ws.Range("D13").Value = "1.639.98"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "4.09"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "65.53"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "26.875.98"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "231.13"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "4.45"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "2.22"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D24").Value = "9.20"
$ws.Range("D25").Value = "145.66"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "15.86"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").Value = "1.453.60"
$ws.Range("E33").Value = "  -4.75%  "
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("E35").Value = "  +5.42%  "
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "0.899"
$ws.Range("E37").Value = "  +7.48%  "
$ws.Range("D38").Value = "0.567"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "65.63"
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "0.969"
$ws.Range("E44").Value = "  +6.69%  "
$ws.Range("D45").Value = "1.809.32"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").Value = "0.778"
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("D47").Value = "90.59"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").Value = "  +3.54%  "
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.57"
$ws.Range("E51").Value = "  +0.65%  "
